$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 / Row 25 swap (Kaspa <-> InternetComputer(DFINITY)) ---
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").Value = "0.163"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "8.55"
$ws.Range("E25").Value = "  +7.68%  "

# --- Per-row Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "56.396.27"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "2.310.22"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "510.50"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").Value = "130.70"
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "0.0998"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").Value = "2.725.53"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "23.38"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "56.378.02"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "2.321.15"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "10.31"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "325.77"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "61.21"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "167.34"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "1.66"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "0.0₃0713"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "18.23"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("D37").Value = "0.879"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").Value = "38.50"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "149.45"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").Value = "3.55"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "273.31"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").Value = "0.0924"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "0.0492"
$ws.Range("D47").Value = "0.550"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "18.05"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "16.90"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "10.99"
$ws.Range("E51").Value = "  +0.38%  "
